$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "30.551.77"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -0.35%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.877.91"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -0.71%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "236.21"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -3.57%  "

$ws.Range("E6").Value = "  +0.10%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.4876"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -1.71%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.2894"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -2.36%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.06669"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -2.19%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "1.872.73"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -0.88%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "16.58"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -2.89%  "

$ws.Range("E12").Value = "  -0.96%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "88.66"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -2.39%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "4.999"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -1.59%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.6510"
$cell.Style = "Normal"

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "30.507.71"
$cell.Style = "Normal"

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.000007842"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -1.33%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +0.12%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "12.99"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -1.71%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "2.115.44"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -0.86%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -0.08%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "4.715"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -2.96%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "194.22"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +10.15%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "6.129"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +0.95%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "9.357"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +0.90%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "156.67"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +1.36%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "18.50"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -1.31%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "1.826"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -5.35%  "

$ws.Range("E29").Value = "  +1.39%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "4.252"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -1.99%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "0.09019"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +1.03%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "3.921"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -2.74%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.05101"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -2.42%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.7210"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -2.82%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.076"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -5.40%  "

$ws.Range("E36").Value = "  +0.75%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.01810"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -3.62%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "2.660"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -1.64%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.9195"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -1.74%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "2.040"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -6.03%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.4382"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +0.16%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "104.56"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -0.97%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.9956"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -0.49%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "5.720"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -1.51%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.1326"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -2.22%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "7.326"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -4.56%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.4015"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +3.39%  "

$ws.Range("E48").Value = "  -0.39%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "8.649"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +1.79%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "1.402"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +1.39%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "33.06"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -0.96%  "
